$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 2000-2009 rows (rows 2-11); this shifts the 2010-2020 data
# (previously rows 12-22) up to rows 2-12.
$ws.Range("A2:A11").EntireRow.Delete()

# Copy the formatting of the last existing data row (2020年, now row 12)
# onto the new row 13 before filling in the 2021年 figures.
$ws.Range("A12").Copy($ws.Range("A13"))

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 426.7239
$ws.Range("C13").Value = 167.2485
$ws.Range("D13").Value = 29.8255
$ws.Range("E13").Value = 2492
$ws.Range("F13").Value = 108.7036
